$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1 + 2: Paragraph that used to read
#   "7. What is Baylor's marching band called? "
# becomes four runs:
#   "7. What is " | "the name of " | "Baylor's marching band " |
#   <bookmarkStart/End name="_GoBack"/> | "? "
# (the _GoBack bookmark used to sit in question 3; adding it here with
# the same name moves it, which also cleans up the old location.)
# ---------------------------------------------------------------------

$marchingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "*marching band called*") {
        $marchingPara = $cand
        break
    }
}

$full = $marchingPara.Range
$s = $full.Start

# "called? " -> "? "  (offsets 34..42 within the paragraph text)
$rCalled = $d.Range($s + 34, $s + 42)
$rCalled.Text = "? "

# Insert "the name of " right after "7. What is " (offset 11)
$insertPoint = $d.Range($s + 11, $s + 11)
$insertPoint.InsertAfter("the name of ")

# Move/re-create the _GoBack bookmark right before the new "? " run
# (now at offset 46: "7. What is the name of Baylor's marching band " = 46 chars)
$bmPoint = $d.Range($s + 46, $s + 46)
$d.Bookmarks.Add("_GoBack", $bmPoint)

# Force a clean run split between "the name of " and "Baylor's marching band "
$split1 = $d.Range($s + 23, $s + 23)
$d.Bookmarks.Add("TempSplit1", $split1)
$d.Bookmarks("TempSplit1").Delete()

# Force a clean run split between "7. What is " and "the name of "
$split2 = $d.Range($s + 11, $s + 11)
$d.Bookmarks.Add("TempSplit2", $split2)
$d.Bookmarks("TempSplit2").Delete()

# ---------------------------------------------------------------------
# Change 3: paragraph "9. What are Baylor's official colors" gains a
# trailing "?" run.
# ---------------------------------------------------------------------

$colorsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "*official colors*") {
        $colorsPara = $cand
        break
    }
}

$cFull = $colorsPara.Range
$cs = $cFull.Start

$cFull.InsertAfter("?")

# Split the new "?" off into its own run
$splitQ9 = $d.Range($cs + 36, $cs + 36)
$d.Bookmarks.Add("TempSplitQ9", $splitQ9)
$d.Bookmarks("TempSplitQ9").Delete()
